# Add data for 2025-11-19
# Updates 2025 (column L) running totals and a couple of prior-year (G/H)
# reclassification corrections across the citywide, by-neighborhood, and
# individual neighborhood breakdown sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5911
$ws.Range("L3").Value = 6438
$ws.Range("G4").Value = 1510
$ws.Range("H4").Value = 1766
$ws.Range("L4").Value = 1585
$ws.Range("L6").Value = 5300
$ws.Range("G7").Value = 24737
$ws.Range("H7").Value = 26082
$ws.Range("L7").Value = 19617

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 171
$ws.Range("L5").Value = 71
$ws.Range("G7").Value = 709
$ws.Range("H7").Value = 810
$ws.Range("L7").Value = 634
$ws.Range("L8").Value = 1290
$ws.Range("L9").Value = 111
$ws.Range("L11").Value = 326
$ws.Range("L27").Value = 174
$ws.Range("L29").Value = 1102
$ws.Range("L33").Value = 886
$ws.Range("L37").Value = 745
$ws.Range("L42").Value = 630
$ws.Range("L43").Value = 146
$ws.Range("L47").Value = 136
$ws.Range("L51").Value = 250
$ws.Range("L52").Value = 409
$ws.Range("L53").Value = 217
$ws.Range("L55").Value = 203
$ws.Range("L63").Value = 56
$ws.Range("L64").Value = 125
$ws.Range("L65").Value = 383
$ws.Range("L66").Value = 59
$ws.Range("L67").Value = 677
$ws.Range("L69").Value = 44
$ws.Range("L72").Value = 81
$ws.Range("L76").Value = 299
$ws.Range("L78").Value = 252
$ws.Range("L81").Value = 17
$ws.Range("L85").Value = 972
$ws.Range("L86").Value = 128
$ws.Range("L91").Value = 264
$ws.Range("L94").Value = 243
$ws.Range("L95").Value = 277
$ws.Range("L98").Value = 104
$ws.Range("L99").Value = 340
$ws.Range("G101").Value = 24737
$ws.Range("H101").Value = 26082
$ws.Range("L101").Value = 19617

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L4").Value = 23
$ws.Range("L6").Value = 66

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 216
$ws.Range("G4").Value = 41
$ws.Range("H4").Value = 32
$ws.Range("L6").Value = 153
$ws.Range("G7").Value = 709
$ws.Range("H7").Value = 810
$ws.Range("L7").Value = 634

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L4").Value = 26
$ws.Range("L7").Value = 326

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 401
$ws.Range("L6").Value = 203
$ws.Range("L7").Value = 972

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L6").Value = 113
$ws.Range("L7").Value = 409

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L3").Value = 56
$ws.Range("L6").Value = 73
$ws.Range("L7").Value = 217

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 385
$ws.Range("L6").Value = 316
$ws.Range("L7").Value = 1290

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 256
$ws.Range("L7").Value = 886

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L4").Value = 17
$ws.Range("L6").Value = 64
$ws.Range("L7").Value = 277

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L4").Value = 39
$ws.Range("L7").Value = 745

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 140
$ws.Range("L3").Value = 125
$ws.Range("L7").Value = 383

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 73
$ws.Range("L7").Value = 340

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L6").Value = 158
$ws.Range("L7").Value = 677

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 429
$ws.Range("L4").Value = 60
$ws.Range("L6").Value = 270
$ws.Range("L7").Value = 1102

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 64
$ws.Range("L7").Value = 299

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 175
$ws.Range("L7").Value = 630

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 85
$ws.Range("L6").Value = 71
$ws.Range("L7").Value = 252

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 68
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 89
$ws.Range("L7").Value = 264

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 36
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 80
$ws.Range("L6").Value = 59

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 91
$ws.Range("L7").Value = 243

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 57
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 171

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 50
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 174

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 68
$ws.Range("L7").Value = 128

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 78
$ws.Range("L3").Value = 79
$ws.Range("L7").Value = 250

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L4").Value = 22
$ws.Range("L7").Value = 146

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item("Sauganash,Forest Glen")
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 17

